$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: restyle from the "header" border (4/5) to the "continuation" border (6/7),
#     matching the style already used by rows 3/6/8/12 (thin bottom border, smaller font).
#     Values already in row 10 (B10=123, C10/D10/E10 shared strings) are left untouched.
$ws.Range("A3:E3").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# --- Row 12: brand-new blank "continuation" row, same style family as row 10.
$ws.Range("A3:E3").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

# --- Row 11: brand-new data row (style comes from the column defaults: 4 for A/B, 5 for C/D/E).
$ws.Range("A11").Value = "SCRIPT/T01P01A/um2406.ssb"
$ws.Range("B11").Value = 101

# C11 contains a literal backslash (" It isn\'t over yet! Good luck!"). A direct
# Range.Value assignment through this host doubles literal backslashes, so the
# text is instead written as a formula-literal (which preserves backslashes
# verbatim) and then flattened back down to a plain value in place.
$ws.Range("C11").Formula = "=`" It isn\'t over yet! Good luck!`""
$ws.Range("C11").Copy()
$ws.Range("C11").PasteSpecial(-4163)

$ws.Range("D11").Value = " Это ещё не конец! Удачи!"
$ws.Range("E11").Value = " Üóï åþæ îå ëïîåø! Ôäàœé!"

# Row heights: row 10 keeps its existing custom height; row 11 gets the standard
# "data row" custom height; row 12 (blank) keeps the default row height.
$ws.Rows.Item(11).RowHeight = 43.2

# Selection follows the last-edited cell.
$ws.Range("D11").Select() | Out-Null
